# Apply the MCH132 collection update: add a Series-level metadata row
# (row 2) below the existing header row, matching the data columns:
# A=identifier  B=alternativeIdentifiers  C=title  D=date_s
# E=levelOfDescription  F=extentAndMedium  G=notes  H=file_path

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the new row's font (Calibri, 10pt, theme text color) once on a
# scratch cell inside the already-used range, then copy that formatting
# onto every row-2 cell so we don't blow up the used range / dimension.
$helper = $ws.Range("I2")
$helper.Font.Name = "Calibri"
$helper.Font.Size = 10
$helper.Font.ThemeColor = 1
$helper.Copy()

$rowCells = @("A2", "C2", "D2", "E2", "F2", "G2", "H2")
foreach ($addr in $rowCells) {
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
}
$helper.Clear()

# Populate the new row's values (B2, C2, D2 and H2 stay blank, like the
# rest of the sheet -- only identifier/levelOfDescription/extentAndMedium/
# notes are known for this series).
$ws.Range("A2").Value = "MCH132-1"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 21L | GRAP COUNT NUMER: NONE"
